# ---------------------------------------------------------------------------
# Add 2022-Q3 data: insert a new worksheet "2022-Q3" right after "总计",
# populate it, and add a corresponding summary row at the top of "总计".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # "总计" summary sheet

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newWs.Name = "2022-Q3"

# Fetch the formatting template sheet *by name* (and only after the insert)
# so it keeps pointing at "2022-Q2" rather than drifting to whatever sheet
# now occupies index 2.
$refWs = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row formatting (bold + border + centred, style index "2")
# and the column-A data style from the reference sheet so the new sheet
# matches the look of its siblings.
$refWs.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)

$refWs.Range("A2").Copy()
$newWs.Range("A2:A3").PasteSpecial(-4122)

# Header row text.
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Data row 2.
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "'000586"
$newWs.Range("C2").Value = "景顺长城中小创精选股票"
$newWs.Range("D2").Value = "'2.21"
$newWs.Range("E2").Value = "'93.50"
$newWs.Range("F2").Value = "'7.87"
$newWs.Range("G2").Value = "'0.1739"
$newWs.Range("H2").Value = 5

# Data row 3.
$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "'260115"
$newWs.Range("C3").Value = "景顺长城中小盘混合"
$newWs.Range("D3").Value = "'0.92"
$newWs.Range("E3").Value = "'92.87"
$newWs.Range("F3").Value = "'6.07"
$newWs.Range("G3").Value = "'0.0558"
$newWs.Range("H3").Value = 2

# The apostrophe-prefix forces Excel to store the numeric-looking values as
# text (matching the source file's inlineStr cells); ClearFormats() then
# strips the "quote prefix" flag Excel adds for that, without touching the
# cell's actual (already-applied) border/font style.
$newWs.Range("B2:G3").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Insert a new row 2 into "总计" for the 2022-Q3 summary figures, pushing
#    the existing quarters down by one row.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(2).Insert()

# Re-derive column A (sequential index = row - 2) for every data row, since
# the insert leaves the literal values attached to the rows they came with.
For ($r = 2; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 2
}

# Give the new A2 the same style as the other column-A cells, then clear any
# stray formatting the row insert may have copied into B2:D2.
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("A2").Value = 0
$ws1.Range("B2:D2").ClearFormats()

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.23

Write-Host "2022-Q3 sheet and summary row added."
